# Combine Flow Location and Flow Type into single column
# (Update "Regulated" -> "Regulated Inflow" and "Glen Canyon Dam" -> "Lake Powell"
# on the ExtremeFlows sheet, and leave a C9 selection as the last active cell.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExtremeFlows")

$ws.Range("C4").Value = "Regulated Inflow"
$ws.Range("C5").Value = "Regulated Inflow"
$ws.Range("B6").Value = "Lake Powell"

$ws.Activate()
$ws.Range("C9").Select()
